$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = "57.062.15"
$ws.Range("E2").Value = "  -1.24%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = "3.084.92"
$ws.Range("E3").Value = "  -0.15%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("E4").Value = "  +0.02%  "

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").Value = "'520.51"
$ws.Range("E5").Value = "  -0.81%  "

# Row 6: 'Solana' -> 'Solana'
$ws.Range("D6").Value = "'135.59"
$ws.Range("E6").Value = "  -3.34%  "

# Row 7: 'USDC' -> 'USDC'
$ws.Range("E7").Value = "  -0.02%  "

# Row 8: 'LidoStakedEther' -> 'LidoStakedEther'
$ws.Range("D8").Value = "3.085.76"
$ws.Range("E8").Value = "  -0.03%  "

# Row 9: 'XRP' -> 'XRP'
$ws.Range("D9").Value = "'0.454"
$ws.Range("E9").Value = "  +2.72%  "

# Row 10: 'Toncoin' -> 'Toncoin'
$ws.Range("E10").Value = "  +2.67%  "

# Row 11: 'Dogecoin' -> 'Dogecoin'
$ws.Range("E11").Value = "  -1.21%  "

# Row 12: 'Cardano' -> 'Cardano'
$ws.Range("E12").Value = "  +2.12%  "

# Row 13: 'WrappedliquidstakedEther2.0' -> 'TRON'
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.135"
$ws.Range("E13").Value = "  +1.80%  "

# Row 14: 'TRON' -> 'WrappedliquidstakedEther2.0'
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.617.04"
$ws.Range("E14").Value = "  -0.10%  "

# Row 15: 'Avalanche' -> 'Avalanche'
$ws.Range("D15").Value = "'25.29"
$ws.Range("E15").Value = "  -0.78%  "

# Row 16: 'ShibaInu' -> 'ShibaInu'
$ws.Range("E16").Value = "  -1.73%  "

# Row 17: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D17").Value = "57.154.03"
$ws.Range("E17").Value = "  -1.13%  "

# Row 18: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D18").Value = "3.085.41"
$ws.Range("E18").Value = "  -0.40%  "

# Row 19: 'Polkadot' -> 'Polkadot'
$ws.Range("E19").Value = "  -3.10%  "

# Row 20: 'Chainlink' -> 'Chainlink'
$ws.Range("D20").Value = "'12.45"
$ws.Range("E20").Value = "  -1.69%  "

# Row 21: 'Uniswap' -> 'Uniswap'
$ws.Range("D21").Value = "'7.85"
$ws.Range("E21").Value = "  -1.13%  "

# Row 22: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D22").Value = "'347.68"
$ws.Range("E22").Value = "  +1.63%  "

# Row 23: 'LEO' -> 'LEO'
$ws.Range("D23").Value = "'5.77"
$ws.Range("E23").Value = "  +1.31%  "

# Row 24: 'Dai' -> 'Dai'
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = "  -0.11%  "

# Row 25: 'Litecoin' -> 'Litecoin'
$ws.Range("D25").Value = "'68.16"
$ws.Range("E25").Value = "  +1.26%  "

# Row 26: 'Polygon' -> 'Polygon'
$ws.Range("E26").Value = "  -2.28%  "

# Row 27: 'Kaspa' -> 'Kaspa'
$ws.Range("E27").Value = "  -1.93%  "

# Row 28: 'Binance-PegBSC-USD' -> 'Binance-PegBSC-USD'
$ws.Range("E28").Value = "  +0.10%  "

# Row 29: 'PEPE' -> 'PEPE'
$ws.Range("D29").Value = "0.0₃0863"
$ws.Range("E29").Value = "  -5.33%  "

# Row 30: 'USDe' -> 'USDe'
$ws.Range("E30").Value = "  -0.09%  "

# Row 31: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D31").Value = "'7.29"
$ws.Range("E31").Value = "  +0.84%  "

# Row 32: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("E32").Value = "  +0.00%  "

# Row 33: 'RenderToken' -> 'RenderToken'
$ws.Range("E33").Value = "  -8.07%  "

# Row 34: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D34").Value = "'20.82"
$ws.Range("E34").Value = "  -0.45%  "

# Row 35: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range("D35").Value = "'4.89"
$ws.Range("E35").Value = "  +6.18%  "

# Row 36: 'Fetch.AI' -> 'Monero'
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'159.45"
$ws.Range("E36").Value = "  +0.49%  "

# Row 37: 'Monero' -> 'Fetch.AI'
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'1.14"
$ws.Range("E37").Value = "  -3.22%  "

# Row 38: 'Aptos' -> 'Aptos'
$ws.Range("D38").Value = "'6.01"
$ws.Range("E38").Value = "  -1.76%  "

# Row 39: 'EnergySwap' -> 'EnergySwap'
$ws.Range("D39").Value = "'25.62"
$ws.Range("E39").Value = "  -1.36%  "

# Row 40: 'ImmutableX' -> 'ImmutableX'
$ws.Range("E40").Value = "  -0.29%  "

# Row 41: 'Hedera' -> 'Hedera'
$ws.Range("D41").Value = "'0.0655"
$ws.Range("E41").Value = "  -1.59%  "

# Row 42: 'Stacks' -> 'Filecoin'
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.05"
$ws.Range("E42").Value = "  +1.59%  "

# Row 43: 'Filecoin' -> 'Stacks'
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.56"
$ws.Range("E43").Value = "  +2.41%  "

# Row 44: 'Mantle' -> 'Mantle'
$ws.Range("E44").Value = "  +1.31%  "

# Row 45: 'Maker' -> 'Maker'
$ws.Range("D45").Value = "2.389.16"
$ws.Range("E45").Value = "  +5.03%  "

# Row 46: 'OKB' -> 'OKB'
$ws.Range("D46").Value = "'36.61"
$ws.Range("E46").Value = "  -0.76%  "

# Row 47: 'FirstDigitalUSD' -> 'FirstDigitalUSD'
$ws.Range("E47").Value = "  +0.05%  "

# Row 48: 'RenzoRestakedETH' -> 'RenzoRestakedETH'
$ws.Range("D48").Value = "3.125.23"
$ws.Range("E48").Value = "  -0.13%  "

# Row 49: 'VeChain' -> 'VeChain'
$ws.Range("E49").Value = "  +0.58%  "

# Row 50: 'ONDO' -> 'ONDO'
$ws.Range("D50").Value = "'0.957"
$ws.Range("E50").Value = "  -2.66%  "

# Row 51: 'Cosmos' -> 'Cosmos'
$ws.Range("D51").Value = "'5.94"
$ws.Range("E51").Value = "  -2.27%  "
